# Apply the scraped-data refresh to the czech-republic cfl-group-a 2023-2024 sheet.
#  1) A handful of same-day fixture rows had their home/away pairing corrected
#     (the two matches played that day were attributed to the wrong row) -
#     fix by swapping the match-specific columns F:V between the two rows.
#  2) Rows 29-31 (also same day) got a 3-way rotation for the same reason.
#  3) Eight new matches (28-29 Oct 2023) were scraped and appended as rows 98-105.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$cols = @("F","G","H","I","J","K","L","M","N","O","P","Q","R","S","T","U","V")

function Get-RowVals {
    param($ws, $row)
    $vals = @{}
    foreach ($c in $cols) {
        $vals[$c] = $ws.Range("$c$row").Value()
    }
    return $vals
}

function Set-RowVals {
    param($ws, $row, $vals)
    foreach ($c in $cols) {
        $ws.Range("$c$row").Value = $vals[$c]
    }
}

function Swap-Rows {
    param($ws, $rowA, $rowB)
    $a = Get-RowVals $ws $rowA
    $b = Get-RowVals $ws $rowB
    Set-RowVals $ws $rowA $b
    Set-RowVals $ws $rowB $a
}

# --- simple two-row swaps (F:V only, A:E stay put) ---
Swap-Rows $ws 10 11
Swap-Rows $ws 14 15
Swap-Rows $ws 16 17
Swap-Rows $ws 43 44
Swap-Rows $ws 46 47
Swap-Rows $ws 48 49
Swap-Rows $ws 62 63
Swap-Rows $ws 70 71
Swap-Rows $ws 76 77

# --- 3-way rotation: new29 = old30, new30 = old31, new31 = old29 ---
$r29 = Get-RowVals $ws 29
$r30 = Get-RowVals $ws 30
$r31 = Get-RowVals $ws 31
Set-RowVals $ws 29 $r30
Set-RowVals $ws 30 $r31
Set-RowVals $ws 31 $r29

# --- append 8 new match rows (98-105), copying formatting from the last row ---
$lastRow = 97
$newCount = 8
$ws.Range("A$lastRow`:V$lastRow").Copy()
$ws.Range("A98:V105").PasteSpecial(-4122)

$newRows = @(
    @{ A=97;  E=45227.42708333334; F="Motorlet Prague";    G=2; H="Vltavin";           I=2;
       J=2.2;  K="27/10/2023 22:13"; L=2.09; M="28/10/2023 09:54";
       N=3.48; O="27/10/2023 22:13"; P=3.48; Q="28/10/2023 09:54";
       R=2.79; S="27/10/2023 22:13"; T=3.12; U="28/10/2023 09:54";
       V="https://www.betexplorer.com/football/czech-republic/cfl-group-a/motorlet-prague-loko-vltavin/t2rnYmFP/" },

    @{ A=98;  E=45227.42708333334; F="Pisek";               G=1; H="Admira Prague";     I=1;
       J=1.79; K="27/10/2023 22:13"; L=2.63; M="28/10/2023 10:04";
       N=3.78; O="27/10/2023 22:13"; P=3.78; Q="28/10/2023 10:04";
       R=3.58; S="27/10/2023 22:13"; T=2.27; U="28/10/2023 10:04";
       V="https://www.betexplorer.com/football/czech-republic/cfl-group-a/pisek-admira-prague/hnk494xP/" },

    @{ A=99;  E=45227.60416666666; F="Domazlice";           G=3; H="Slavia Prague B";   I=2;
       J=2.06; K="28/10/2023 13:42"; L=2.06; M="28/10/2023 14:15";
       N=3.55; O="28/10/2023 13:42"; P=3.5;  Q="28/10/2023 14:15";
       R=3.14; S="28/10/2023 13:42"; T=3.18; U="28/10/2023 14:15";
       V="https://www.betexplorer.com/football/czech-republic/cfl-group-a/domazlice-slavia-prague/QZJISRxm/" },

    @{ A=100; E=45227.75;           F="Karlovy Vary";       G=1; H="FK Robstav";        I=1;
       J=2.7;  K="28/10/2023 13:45"; L=2.88; M="28/10/2023 17:49";
       N=3.44; O="28/10/2023 13:45"; P=3.42; Q="28/10/2023 17:49";
       R=2.32; S="28/10/2023 13:45"; T=2.25; U="28/10/2023 17:49";
       V="https://www.betexplorer.com/football/czech-republic/cfl-group-a/karlovy-vary-fk-robstav/vwOMRohg/" },

    @{ A=101; E=45228.42708333334; F="Dukla Prague B";      G=3; H="Povltavska FA";     I=3;
       J=2.29; K="29/10/2023 01:42"; L=2.15; M="29/10/2023 10:14";
       N=3.49; O="29/10/2023 01:42"; P=3.61; Q="29/10/2023 10:05";
       R=2.65; S="29/10/2023 01:42"; T=2.91; U="29/10/2023 10:14";
       V="https://www.betexplorer.com/football/czech-republic/cfl-group-a/dukla-prague-povltavska-fa/vDhROqxC/" },

    @{ A=102; E=45228.42708333334; F="Bohemians 1905 B";    G=0; H="Kraluv Dvur";       I=0;
       J=1.5;  K="29/10/2023 01:42"; L=1.7;  M="29/10/2023 10:07";
       N=4.57; O="29/10/2023 01:42"; P=4.07; Q="29/10/2023 10:07";
       R=4.59; S="29/10/2023 01:42"; T=3.96; U="29/10/2023 10:06";
       V="https://www.betexplorer.com/football/czech-republic/cfl-group-a/bohemians-1905-kraluv-dvur/YaNQQ57a/" },

    @{ A=103; E=45228.4375;         F="Plzen B";             G=2; H="Ceske Budejovice B"; I=1;
       J=2.17; K="29/10/2023 01:42"; L=1.93; M="29/10/2023 10:29";
       N=3.64; O="29/10/2023 01:42"; P=3.91; Q="29/10/2023 10:28";
       R=2.81; S="29/10/2023 01:42"; T=3.22; U="29/10/2023 10:29";
       V="https://www.betexplorer.com/football/czech-republic/cfl-group-a/plzen-ceske-budejovice/rBUDT7Nt/" },

    @{ A=104; E=45228.60416666666; F="Taborsko akademie";   G=1; H="Hostoun";           I=0;
       J=2.6;  K="29/10/2023 10:25"; L=2.25; M="29/10/2023 14:29";
       N=3.4;  O="29/10/2023 10:25"; P=3.43; Q="29/10/2023 14:29";
       R=2.43; S="29/10/2023 10:25"; T=2.87; U="29/10/2023 14:29";
       V="https://www.betexplorer.com/football/czech-republic/cfl-group-a/taborsko-akademie-hostoun/CjMUPPM5/" }
)

$rowNum = 98
foreach ($data in $newRows) {
    $ws.Range("A$rowNum").Value = $data.A
    $ws.Range("B$rowNum").Value = "czech-republic"
    $ws.Range("C$rowNum").Value = "cfl-group-a"
    $ws.Range("D$rowNum").Value = "2023-2024"
    $ws.Range("E$rowNum").Value = $data.E
    foreach ($c in $cols) {
        $ws.Range("$c$rowNum").Value = $data[$c]
    }
    $rowNum++
}

Write-Output "applied czech-republic_cfl-group-a_2023-2024 update"
